$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1431.4286
$ws.Range("I135").Value = 984.1
$ws.Range("K135").Value = 8856.9
$ws.Range("M135").Value = -6321.9
$ws.Range("H138").Value = 4589.613
$ws.Range("J138").Value = 4865.8213
$ws.Range("L138").Value = 14597.4639
$ws.Range("N138").Value = -24877.4639

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 391.33334
$ws.Range("I4").Value = 87
$ws.Range("K4").Value = 87
$ws.Range("M4").Value = 29
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 60
$ws.Range("K5").Value = 60
$ws.Range("M5").Value = 52
$ws.Range("H32").Value = 16688.945
$ws.Range("I32").Value = 7516.871
$ws.Range("J32").Value = 29051.305
$ws.Range("K32").Value = 7516.871
$ws.Range("L32").Value = 29051.305
$ws.Range("M32").Value = -7229.871
$ws.Range("N32").Value = -29625.305
$ws.Range("H63").Value = 5696.4375
$ws.Range("I63").Value = 4542.5713
$ws.Range("J63").Value = 6593.8887
$ws.Range("K63").Value = 4542.5713
$ws.Range("L63").Value = 6593.8887
$ws.Range("M63").Value = -3856.5713
$ws.Range("N63").Value = -7965.8887
$ws.Range("H66").Value = 5696.4375
$ws.Range("I66").Value = 4542.5713
$ws.Range("J66").Value = 6593.8887
$ws.Range("K66").Value = 22712.8565
$ws.Range("L66").Value = 32969.4435
$ws.Range("M66").Value = -19280.8565
$ws.Range("N66").Value = -39833.4435
$ws.Range("H122").Value = 718135.6
$ws.Range("I122").Value = 1669331.1
$ws.Range("J122").Value = 4739
$ws.Range("K122").Value = 5007993.300000001
$ws.Range("L122").Value = 14217
$ws.Range("M122").Value = -5005543.300000001
$ws.Range("N122").Value = -19117
$ws.Range("H126").Value = 7989
$ws.Range("I126").Value = 7989
$ws.Range("K126").Value = 23967
$ws.Range("M126").Value = -21497
$ws.Range("H132").Value = 2060.16
$ws.Range("I132").Value = 2000.5454
$ws.Range("K132").Value = 6001.6362
$ws.Range("M132").Value = -3471.6362

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 60
$ws.Range("K4").Value = 60
$ws.Range("M4").Value = 55
$ws.Range("H43").Value = 150000
$ws.Range("J43").Value = 150000
$ws.Range("L43").Value = 150000
$ws.Range("N43").Value = -150362
$ws.Range("H100").Value = 32242
$ws.Range("J100").Value = 32242
$ws.Range("L100").Value = 32242
$ws.Range("N100").Value = -34406
$ws.Range("H105").Value = 3491.2432
$ws.Range("I105").Value = 2815.12
$ws.Range("K105").Value = 2815.12
$ws.Range("M105").Value = -1068.12
$ws.Range("H132").Value = 125000
$ws.Range("J132").Value = 125000
$ws.Range("L132").Value = 125000
$ws.Range("N132").Value = -135120
$ws.Range("H134").Value = 5333
$ws.Range("I134").Value = 5299.5
$ws.Range("K134").Value = 15898.5
$ws.Range("M134").Value = -13363.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 375.7143
$ws.Range("I22").Value = 399.5
$ws.Range("J22").Value = 344
$ws.Range("K22").Value = 399.5
$ws.Range("L22").Value = 344
$ws.Range("M22").Value = -49.5
$ws.Range("N22").Value = -1044
$ws.Range("H28").Value = 16854
$ws.Range("J28").Value = 16854
$ws.Range("L28").Value = 16854
$ws.Range("N28").Value = -17344
$ws.Range("H58").Value = 3395.6316
$ws.Range("I58").Value = 2237.4443
$ws.Range("J58").Value = 4438
$ws.Range("K58").Value = 2237.4443
$ws.Range("L58").Value = 4438
$ws.Range("M58").Value = -2034.4443
$ws.Range("N58").Value = -4844
$ws.Range("H103").Value = 13456.125
$ws.Range("I103").Value = 11092.714
$ws.Range("K103").Value = 11092.714
$ws.Range("M103").Value = -9920.714
$ws.Range("H136").Value = 3395.6316
$ws.Range("I136").Value = 2237.4443
$ws.Range("J136").Value = 4438
$ws.Range("K136").Value = 6712.3329
$ws.Range("L136").Value = 13314
$ws.Range("M136").Value = -4162.3329
$ws.Range("N136").Value = -18414
$ws.Range("H141").Value = 155800
$ws.Range("J141").Value = 155800
$ws.Range("L141").Value = 155800
$ws.Range("N141").Value = -166160

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 299.75
$ws.Range("I97").Value = 99.5
$ws.Range("K97").Value = 298.5
$ws.Range("M97").Value = 197.5
$ws.Range("H120").Value = 15495
$ws.Range("I120").Value = 7000
$ws.Range("J120").Value = 15942.105
$ws.Range("K120").Value = 21000
$ws.Range("L120").Value = 47826.315
$ws.Range("M120").Value = -16162
$ws.Range("N120").Value = -57502.315

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 319.33334
$ws.Range("I2").Value = 99.625
$ws.Range("K2").Value = 99.625
$ws.Range("M2").Value = 13.375
$ws.Range("H62").Value = 60000
$ws.Range("J62").Value = 60000
$ws.Range("L62").Value = 60000
$ws.Range("N62").Value = -61372
$ws.Range("H65").Value = 60000
$ws.Range("J65").Value = 60000
$ws.Range("L65").Value = 180000
$ws.Range("N65").Value = -186864
$ws.Range("H98").Value = 31643
$ws.Range("J98").Value = 31643
$ws.Range("L98").Value = 31643
$ws.Range("N98").Value = -37633
$ws.Range("H101").Value = 4999.6665
$ws.Range("J101").Value = 4999.6665
$ws.Range("L101").Value = 4999.6665
$ws.Range("N101").Value = -11489.6665
$ws.Range("H132").Value = 8546.75
$ws.Range("I132").Value = 8433
$ws.Range("K132").Value = 25299
$ws.Range("M132").Value = -22769

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2622.3333
$ws.Range("J46").Value = 3354.0908
$ws.Range("L46").Value = 3354.0908
$ws.Range("N46").Value = -3730.0908
$ws.Range("H55").Value = 384.55
$ws.Range("I55").Value = 330.6
$ws.Range("K55").Value = 330.6
$ws.Range("M55").Value = -157.6

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 57868.223
$ws.Range("I14").Value = 73200.57000000001
$ws.Range("J14").Value = 4205
$ws.Range("K14").Value = 73200.57000000001
$ws.Range("L14").Value = 4205
$ws.Range("M14").Value = -73032.57000000001
$ws.Range("N14").Value = -4541
$ws.Range("H64").Value = 68000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 68000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 68000
$ws.Range("M64").Value = ""
$ws.Range("N64").Value = -68496
$ws.Range("H67").Value = 68000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 68000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 68000
$ws.Range("M67").Value = ""
$ws.Range("N67").Value = -69716
$ws.Range("H113").Value = 1224
$ws.Range("I113").Value = 957.5
$ws.Range("K113").Value = 2872.5
$ws.Range("M113").Value = -702.5
$ws.Range("H136").Value = 51735.1
$ws.Range("J136").Value = 127219
$ws.Range("L136").Value = 381657
$ws.Range("N136").Value = -386757
